# Re-order a set of rows in "Germany Bundesliga I": each listed row's
# columns B:AD (everything except the running index in column A) get
# replaced by the B:AD values that currently live in a different row,
# per the mapping below. Column A (the sequential counter) is left
# untouched on every row.
#
# Because several of these reassignments form cycles (e.g. row 39 takes
# row 43's data, row 40 takes row 39's old data, ..., row 43 takes row
# 42's old data), we must snapshot every source row's B:AD values BEFORE
# writing any of them back out, then apply all the writes afterward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (source row's current B:AD values become the
# target row's new B:AD values)
$rowMap = @{
    31  = 32
    32  = 31
    39  = 43
    40  = 39
    41  = 40
    42  = 41
    43  = 42
    77  = 78
    78  = 79
    79  = 77
    149 = 150
    150 = 149
    173 = 174
    174 = 177
    175 = 173
    177 = 175
    192 = 194
    193 = 192
    194 = 193
    202 = 204
    203 = 202
    204 = 203
    238 = 240
    240 = 238
    301 = 303
    302 = 301
    303 = 302
}

# Snapshot B:AD for every row that participates (as source or target).
$snapshot = @{}
foreach ($target in $rowMap.Keys) {
    $source = $rowMap[$target]
    if (-not $snapshot.ContainsKey($source)) {
        $snapshot[$source] = $ws.Range("B$($source):AD$($source)").Value2
    }
}

# Now write the snapshotted values into their new homes.
foreach ($target in $rowMap.Keys) {
    $source = $rowMap[$target]
    $ws.Range("B$($target):AD$($target)").Value2 = $snapshot[$source]
}
